# Bibliografia paragraph: turn the single run of concatenated references
# into separate runs divided by manual line breaks (<w:br/>), one per
# numbered reference.

$d = $word.ActiveDocument

# 1) break between ref 1 and ref 2
$d.Content.Find.Execute(
    "Editora Blucher, 2021. 2. Cisternas, J.R.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Editora Blucher, 2021. ^l2. Cisternas, J.R.", 2) | Out-Null

# 2) break between ref 2 and ref 3
$d.Content.Find.Execute(
    "São Paulo: Atheneu, 2005. 3. Nelson, D.L.", $true, $false, $false, $false, $false,
    $true, 1, $false, "São Paulo: Atheneu, 2005. ^l3. Nelson, D.L.", 2) | Out-Null

# 3) break between ref 3 and ref 4
$d.Content.Find.Execute(
    "Artmed Editora, 2022. 4. Voet, D.,", $true, $false, $false, $false, $false,
    $true, 1, $false, "Artmed Editora, 2022. ^l4. Voet, D.,", 2) | Out-Null

# 4) break between ref 4 and ref 5
$d.Content.Find.Execute(
    "Artmed Editora, 2014. 5. Vitolo, M.,", $true, $false, $false, $false, $false,
    $true, 1, $false, "Artmed Editora, 2014. ^l5. Vitolo, M.,", 2) | Out-Null
